$wb = $excel.ActiveWorkbook

# Rename the first sheet ("Population Types" -> "Population Attributes")
$popSheet = $wb.Worksheets.Item("Population Types")
$popSheet.Name = "Population Attributes"

# Make this sheet the active/selected tab instead of "Characteristics"
$popSheet.Activate()
